$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: 2022-12-25, Nereocystis.luetkeana = "y" (column D)
$ws.Range("A9").Value = 2022
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 25
$ws.Range("D9").Value = "y"

# Row 10: 2023-02-18, Nereocystis.luetkeana = "y" (column D)
$ws.Range("A10").Value = 2023
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 18
$ws.Range("D10").Value = "y"

# Update the selected cell to match the saved view state
$ws.Range("E10").Select()
